# Updated cryptos list — applies the per-cell Price/Volume(1h) refresh
# plus the row 50/51 ranking shift (BabyDogeCoin entering the top list,
# Algorand sliding from rank 50 to rank 51, USDD dropping off).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.005.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "'1.564.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'207.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'22.14"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "'0.0602"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("D11").Value = "'0.0858"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "'1.786.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "'1.564.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "'3.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "'62.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "'27.002.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").Value = "'0.0₃0705"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").Value = "'216.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D20").Value = "'7.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").Value = "'9.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "'153.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").Value = "'6.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "'15.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").Value = "'1.422.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("D35").Value = "'1.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.26%  "
$ws.Range("E36").Value = "  +10.50%  "
$ws.Range("D37").Value = "'2.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.13%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "'0.531"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("D41").Value = "'0.808"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("D44").Value = "'1.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").Value = "'64.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").Value = "'1.700.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").Value = "'87.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").Value = "'0.0520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.0₇0997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0960"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.20%  "
